$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Dummy" header in column F (new shared string + new column of data)
$ws.Range("F1").Value = "Dummy"

# Move/resize the existing chart the same way the author repositioned it
# (captured from the target anchor: from col6/off103907,row0/off69273 to
#  col14/off505689,row17/off51954 translated into points).
$co = $ws.ChartObjects(1)
$co.Left = 445.1787238558071
$co.Top = 5.4545669291338585
$co.Width = 499.1363779527559
$co.Height = 258.73629921259845

# Selection ends up on F2 after the edit
$ws.Range("F2").Select()
